$d = $word.ActiveDocument

# NOTE: Find.Execute() on a Range obtained via Table/Cell .Range (or a
# .Duplicate() of one) ignores that range's own Start and always searches
# from the top of the document story instead of the range's own start.
# Re-building the range explicitly with $d.Range(start, end) right before
# calling .Find makes the search honor the intended [start, end) window, so
# every lookup below goes through that pattern. The found hit is likewise
# re-wrapped in a fresh $d.Range(...) before any write, since writing
# straight back into the range object that ran .Find mistargets the edit.

# ---------------------------------------------------------------------------
# Edit 1: title-page table (Tables(1)), row 1 / col 2 — the cell holding the
# underlined project-title sentence ending "...для ІТ-галузі. API машинного
# навчання". The standalone run containing exactly ". " (right after
# "...для ІТ-галузі") becomes " (комплексна тема). ".
# ---------------------------------------------------------------------------
$cell1 = $d.Tables.Item(1).Cell(1, 2)
$c1Start = $cell1.Range.Start
$c1End = $cell1.Range.End

if ($cell1.Range.Text -notmatch "Вебзастосунок") {
    throw "Edit 1: Tables(1).Cell(1,2) does not look like the expected project-title cell"
}

$scope1 = $d.Range($c1Start, $c1End)
$hit1 = $scope1.Find.Execute(". ", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
if ($hit1 -and $scope1.Start -ge $c1Start -and $scope1.End -le $c1End) {
    $target1 = $d.Range($scope1.Start, $scope1.End)
    $target1.Text = " (комплексна тема). "
} else {
    throw "Edit 1: target run '. ' not found in Tables(1).Cell(1,2)"
}

# ---------------------------------------------------------------------------
# Edit 2: project-info table (Tables(5)), row 1 / col 2 — same underlined
# sentence, repeated in the project-info table. Insert " (комплексна тема)"
# right after "...для ІТ-галузі" and before the following ". API " run.
# ---------------------------------------------------------------------------
$cell2 = $d.Tables.Item(5).Cell(1, 2)
$c2Start = $cell2.Range.Start
$c2End = $cell2.Range.End

if ($cell2.Range.Text -notmatch "Вебзастосунок") {
    throw "Edit 2: Tables(5).Cell(1,2) does not look like the expected project-title cell"
}

$scope2 = $d.Range($c2Start, $c2End)
$hit2 = $scope2.Find.Execute("для ІТ-галузі", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
if ($hit2 -and $scope2.Start -ge $c2Start -and $scope2.End -le $c2End) {
    $insertionPoint = $d.Range($scope2.End, $scope2.End)
    $insertionPoint.InsertAfter(" (комплексна тема)")
} else {
    throw "Edit 2: target text 'для ІТ-галузі' not found in Tables(5).Cell(1,2)"
}

Write-Output "Edit 1 applied: $hit1; Edit 2 applied: $hit2"
